# DNUG Presentation update - BP
# 1) Title slide: "Angular 4" -> "Angular 2+"
# 2) Insert a new slide (position 13) "Things not covered" with a bulleted list
# 3) Remove the "Replace both controllers and directives..." bullet from the
#    "Components" slide

$p = $ppt.ActivePresentation

# --- 1) Update the deck title on the title slide ---------------------------
$titleSlide = $p.Slides.Item(1)
$titleSlide.Shapes.Item(1).TextFrame.TextRange.Text = "Angular 2+"

# --- 2) Insert new "Things not covered" slide before "3rd party libraries" -
# ppLayoutText (=2) matches the "Title and Content" layout used by the
# surrounding slides.
$newSlide = $p.Slides.Add(13, 2)

$newSlide.Shapes.Item(1).TextFrame.TextRange.Text = "Things not covered"

$bodyRange = $newSlide.Shapes.Item(2).TextFrame.TextRange
$bodyRange.Text = "Directives"
[void]$bodyRange.InsertAfter("`rPipes")
[void]$bodyRange.InsertAfter("`rGuards")
[void]$bodyRange.InsertAfter("`rModules")
[void]$bodyRange.InsertAfter("`rUnit tests")
[void]$bodyRange.InsertAfter("`rE2E tests")

# --- 3) Drop the redundant "Angular 1" bullet from the Components slide ----
$componentsSlide = $p.Slides.Item(7)
$componentsBody = $componentsSlide.Shapes.Item(2).TextFrame.TextRange
[void]$componentsBody.Paragraphs(5).Delete()
